# The deck's design theme is switched from the custom "Integral" theme
# back to the default Office "Office Theme" palette (Design > Themes >
# "Office" in the PowerPoint UI). This repaints the 12 theme colors used
# by the slide master / layouts / slides (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) with the stock Office theme RGB values.

function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# Office theme ("Office Theme") color scheme, in clrScheme document order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$cs.Colors(1).RGB  = RGB 0   0   0     # dk1      000000
$cs.Colors(2).RGB  = RGB 255 255 255   # lt1      FFFFFF
$cs.Colors(3).RGB  = RGB 68  84  106   # dk2      44546A
$cs.Colors(4).RGB  = RGB 231 230 230   # lt2      E7E6E6
$cs.Colors(5).RGB  = RGB 91  155 213   # accent1  5B9BD5
$cs.Colors(6).RGB  = RGB 237 125 49    # accent2  ED7D31
$cs.Colors(7).RGB  = RGB 165 165 165   # accent3  A5A5A5
$cs.Colors(8).RGB  = RGB 255 192 0     # accent4  FFC000
$cs.Colors(9).RGB  = RGB 68  114 196   # accent5  4472C4
$cs.Colors(10).RGB = RGB 112 173 71    # accent6  70AD47
$cs.Colors(11).RGB = RGB 5   99  193   # hlink    0563C1
$cs.Colors(12).RGB = RGB 149 79  114   # folHlink 954F72
